$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix letter case inconsistency: "Vocabulary Code" -> "Vocabulary code"
$ws.Range("H4").Value = "Vocabulary code"

# Update selection to reflect where the edit was made
$ws.Range("H4").Select()
